# Minor improvements mainly in Documentation
# - Reword several "Test Case Objective" cells (rows 23-33) to clarify
#   whether the scenario is performed by a "Visitor" or a logged-in "User".
# - Update the saved sheet view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Test Case Objective wording (column E, rows 23-33) ---
$ws.Range("E23").Value = "Visitor Navigate to Login page"
$ws.Range("E24").Value = "Visitor Navigate to Register page"
$ws.Range("E25").Value = "Visitor Navigate to Home page through Logo link"
$ws.Range("E26").Value = "User Login and Enter Account Management"
$ws.Range("E27").Value = "User Navigate to Create Article and click Cancel"
$ws.Range("E28").Value = "User Login and Read an Article"
$ws.Range("E29").Value = "User Login and Select Article to Delete"
$ws.Range("E30").Value = "User Login and Select Article to Edit"
$ws.Range("E31").Value = "User Login and navigate to Password change Page"
$ws.Range("E32").Value = "User Login and Comment created Article"
$ws.Range("E33").Value = "Visitor View article Author "

# --- Update the saved view: scrolled so row 4 is at the top, with E28 selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E28").Select()
